{"js": "// Replace the date line and each two-digit \u00f7 one-digit practice answer\n// in the table with the new values from the latest generated output.\nconst replacements = [\n  [\"2025-02-28 Friday\", \"2025-03-01 Saturday\"],\n  [\"42\u00f75=8, 2\", \"64\u00f73=21, 1\"],\n  [\"42\u00f79=4, 6\", \"25\u00f73=8, 1\"],\n  [\"45\u00f76=7, 3\", \"62\u00f73=20, 2\"],\n  [\"88\u00f74=22, 0\", \"99\u00f73=33, 0\"],\n  [\"69\u00f78=8, 5\", \"30\u00f77=4, 2\"],\n  [\"24\u00f77=3, 3\", \"84\u00f76=14, 0\"],\n  [\"67\u00f79=7, 4\", \"37\u00f78=4, 5\"],\n  [\"21\u00f78=2, 5\", \"65\u00f72=32, 1\"],\n  [\"52\u00f73=17, 1\", \"37\u00f77=5, 2\"],\n  [\"55\u00f77=7, 6\", \"98\u00f73=32, 2\"],\n  [\"45\u00f75=9, 0\", \"77\u00f74=19, 1\"],\n  [\"33\u00f77=4, 5\", \"70\u00f73=23, 1\"],\n  [\"85\u00f74=21, 1\", \"64\u00f72=32, 0\"],\n  [\"89\u00f76=14, 5\", \"46\u00f74=11, 2\"],\n  [\"59\u00f72=29, 1\", \"44\u00f75=8, 4\"],\n  [\"85\u00f78=10, 5\", \"53\u00f77=7, 4\"],\n  [\"38\u00f75=7, 3\", \"49\u00f72=24, 1\"],\n  [\"63\u00f78=7, 7\", \"46\u00f77=6, 4\"],\n  [\"41\u00f79=4, 5\", \"39\u00f77=5, 4\"],\n  [\"44\u00f72=22, 0\", \"13\u00f73=4, 1\"],\n  [\"41\u00f73=13, 2\", \"10\u00f78=1, 2\"],\n  [\"21\u00f76=3, 3\", \"32\u00f78=4, 0\"],\n  [\"69\u00f73=23, 0\", \"46\u00f75=9, 1\"],\n  [\"68\u00f73=22, 2\", \"28\u00f79=3, 1\"],\n  [\"82\u00f78=10, 2\", \"89\u00f79=9, 8\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const item of found.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update master to output generated at 503736d\n# Replace the date line and each two-digit / one-digit division answer\n# in the practice table with the newly generated values.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-02-28 Friday\", \"2025-03-01 Saturday\"),\n    @(\"42\u00f75=8, 2\", \"64\u00f73=21, 1\"),\n    @(\"42\u00f79=4, 6\", \"25\u00f73=8, 1\"),\n    @(\"45\u00f76=7, 3\", \"62\u00f73=20, 2\"),\n    @(\"88\u00f74=22, 0\", \"99\u00f73=33, 0\"),\n    @(\"69\u00f78=8, 5\", \"30\u00f77=4, 2\"),\n    @(\"24\u00f77=3, 3\", \"84\u00f76=14, 0\"),\n    @(\"67\u00f79=7, 4\", \"37\u00f78=4, 5\"),\n    @(\"21\u00f78=2, 5\", \"65\u00f72=32, 1\"),\n    @(\"52\u00f73=17, 1\", \"37\u00f77=5, 2\"),\n    @(\"55\u00f77=7, 6\", \"98\u00f73=32, 2\"),\n    @(\"45\u00f75=9, 0\", \"77\u00f74=19, 1\"),\n    @(\"33\u00f77=4, 5\", \"70\u00f73=23, 1\"),\n    @(\"85\u00f74=21, 1\", \"64\u00f72=32, 0\"),\n    @(\"89\u00f76=14, 5\", \"46\u00f74=11, 2\"),\n    @(\"59\u00f72=29, 1\", \"44\u00f75=8, 4\"),\n    @(\"85\u00f78=10, 5\", \"53\u00f77=7, 4\"),\n    @(\"38\u00f75=7, 3\", \"49\u00f72=24, 1\"),\n    @(\"63\u00f78=7, 7\", \"46\u00f77=6, 4\"),\n    @(\"41\u00f79=4, 5\", \"39\u00f77=5, 4\"),\n    @(\"44\u00f72=22, 0\", \"13\u00f73=4, 1\"),\n    @(\"41\u00f73=13, 2\", \"10\u00f78=1, 2\"),\n    @(\"21\u00f76=3, 3\", \"32\u00f78=4, 0\"),\n    @(\"69\u00f73=23, 0\", \"46\u00f75=9, 1\"),\n    @(\"68\u00f73=22, 2\", \"28\u00f79=3, 1\"),\n    @(\"82\u00f78=10, 2\", \"89\u00f79=9, 8\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
